$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41 (shifts existing rows 41..149 down to 42..150)
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new weekly record
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 45028
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = 100112035
$ws.Range("G41").Value = "Bruselas (repollito)"
$ws.Range("H41").Value = "Sin especificar"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 30
$ws.Range("K41").Value = 30000
$ws.Range("L41").Value = 30000
$ws.Range("M41").Value = 30000
$ws.Range("N41").Value = "$/malla 15 kilos"
$ws.Range("O41").Value = "Región Metropolitana"
$ws.Range("P41").Value = 2000
$ws.Range("Q41").Value = 15
$ws.Range("R41").Value = "Hortaliza"
